$d = $word.ActiveDocument

# Update the title/date paragraph.
# (Note: this runtime's Find.Execute does not actually confine its search to
# the Range/Find object it was invoked on -- it always matches against the
# whole document. Directly assigning to a Range's .Text property, however,
# is correctly scoped to that Range, so it is used everywhere here to avoid
# accidentally touching one of the several duplicated cell strings below.)
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "2025-02-01 Saturday"

# Update the table of division problems, cell by cell, addressed by
# (row, column) position so duplicated source strings elsewhere in the
# table are never accidentally matched/replaced.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "89÷5=17, 4" },
    @{ Row = 1;  Col = 2; New = "41÷3=13, 2" },
    @{ Row = 1;  Col = 3; New = "31÷4=7, 3" },
    @{ Row = 1;  Col = 4; New = "65÷3=21, 2" },
    @{ Row = 1;  Col = 5; New = "31÷5=6, 1" },

    @{ Row = 5;  Col = 1; New = "43÷5=8, 3" },
    @{ Row = 5;  Col = 2; New = "87÷4=21, 3" },
    @{ Row = 5;  Col = 3; New = "79÷5=15, 4" },
    @{ Row = 5;  Col = 4; New = "43÷2=21, 1" },
    @{ Row = 5;  Col = 5; New = "15÷5=3, 0" },

    @{ Row = 9;  Col = 1; New = "17÷8=2, 1" },
    @{ Row = 9;  Col = 2; New = "95÷4=23, 3" },
    @{ Row = 9;  Col = 3; New = "67÷6=11, 1" },
    @{ Row = 9;  Col = 4; New = "95÷6=15, 5" },
    @{ Row = 9;  Col = 5; New = "94÷6=15, 4" },

    @{ Row = 13; Col = 1; New = "46÷6=7, 4" },
    @{ Row = 13; Col = 2; New = "41÷2=20, 1" },
    @{ Row = 13; Col = 3; New = "31÷8=3, 7" },
    @{ Row = 13; Col = 4; New = "61÷5=12, 1" },
    @{ Row = 13; Col = 5; New = "65÷6=10, 5" },

    @{ Row = 17; Col = 1; New = "55÷8=6, 7" },
    @{ Row = 17; Col = 2; New = "47÷8=5, 7" },
    @{ Row = 17; Col = 3; New = "79÷2=39, 1" },
    @{ Row = 17; Col = 4; New = "90÷2=45, 0" },
    @{ Row = 17; Col = 5; New = "56÷7=8, 0" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
